$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '67.511.64'
Set-TextValue $ws.Range('E2') '  -3.20%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.626.34'
Set-TextValue $ws.Range('E3') '  -3.70%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.36%  '

# Row 5
Set-TextValue $ws.Range('D5') '589.36'
Set-TextValue $ws.Range('E5') '  -2.05%  '

# Row 6
Set-TextValue $ws.Range('D6') '185.26'
Set-TextValue $ws.Range('E6') '  -1.01%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.614'
Set-TextValue $ws.Range('E7') '  -3.81%  '

# Row 8
Set-TextValue $ws.Range('E8') '  +0.45%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.678'
Set-TextValue $ws.Range('E9') '  -7.29%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.147'
Set-TextValue $ws.Range('E10') '  -11.22%  '

# Row 11
Set-TextValue $ws.Range('D11') '54.20'
Set-TextValue $ws.Range('E11') '  -5.87%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0000255'
Set-TextValue $ws.Range('E12') '  -14.36%  '

# Row 13
Set-TextValue $ws.Range('D13') '10.02'
Set-TextValue $ws.Range('E13') '  -8.05%  '

# Row 14
Set-TextValue $ws.Range('D14') '4.202.10'
Set-TextValue $ws.Range('E14') '  -4.08%  '

# Row 15
Set-TextValue $ws.Range('D15') '3.629.09'
Set-TextValue $ws.Range('E15') '  -3.98%  '

# Row 16
Set-TextValue $ws.Range('E16') '  -0.51%  '

# Row 17
Set-TextValue $ws.Range('B17') 'Chainlink'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D17') '18.48'
Set-TextValue $ws.Range('E17') '  -6.00%  '

# Row 18
Set-TextValue $ws.Range('B18') 'WrappedBTC'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D18') '67.288.25'
Set-TextValue $ws.Range('E18') '  -3.35%  '

# Row 19
Set-TextValue $ws.Range('D19') '12.33'
Set-TextValue $ws.Range('E19') '  -5.71%  '

# Row 20
Set-TextValue $ws.Range('E20') '  -6.03%  '

# Row 21
Set-TextValue $ws.Range('D21') '400.67'
Set-TextValue $ws.Range('E21') '  -4.03%  '

# Row 22
Set-TextValue $ws.Range('D22') '4.33'
Set-TextValue $ws.Range('E22') '  -7.26%  '

# Row 23
Set-TextValue $ws.Range('D23') '85.91'
Set-TextValue $ws.Range('E23') '  -4.87%  '

# Row 24
Set-TextValue $ws.Range('D24') '2.87'
Set-TextValue $ws.Range('E24') '  -7.17%  '

# Row 25
Set-TextValue $ws.Range('D25') '12.40'
Set-TextValue $ws.Range('E25') '  -5.22%  '

# Row 26
Set-TextValue $ws.Range('D26') '6.07'
Set-TextValue $ws.Range('E26') '  -0.37%  '

# Row 27
Set-TextValue $ws.Range('D27') '10.36'
Set-TextValue $ws.Range('E27') '  -8.34%  '

# Row 28
Set-TextValue $ws.Range('E28') '  -9.90%  '

# Row 29
Set-TextValue $ws.Range('D29') '9.08'
Set-TextValue $ws.Range('E29') '  -5.87%  '

# Row 30
Set-TextValue $ws.Range('D30') '31.38'
Set-TextValue $ws.Range('E30') '  -5.86%  '

# Row 31
Set-TextValue $ws.Range('D31') '6.81'
Set-TextValue $ws.Range('E31') '  -9.22%  '

# Row 32
Set-TextValue $ws.Range('D32') '66.10'
Set-TextValue $ws.Range('E32') '  +1.37%  '

# Row 33
Set-TextValue $ws.Range('D33') '11.94'
Set-TextValue $ws.Range('E33') '  -5.43%  '

# Row 34
Set-TextValue $ws.Range('D34') '595.48'
Set-TextValue $ws.Range('E34') '  -2.56%  '

# Row 35
Set-TextValue $ws.Range('E35') '  -5.88%  '

# Row 36
Set-TextValue $ws.Range('D36') '41.70'
Set-TextValue $ws.Range('E36') '  -6.90%  '

# Row 37
Set-TextValue $ws.Range('D37') '1.00'
Set-TextValue $ws.Range('E37') '  +0.19%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.998'
Set-TextValue $ws.Range('E38') '  -0.69%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.378'
Set-TextValue $ws.Range('E39') '  -7.99%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.0₃0744'
Set-TextValue $ws.Range('E40') '  -18.89%  '

# Row 41
Set-TextValue $ws.Range('E41') '  -4.55%  '

# Row 42
Set-TextValue $ws.Range('D42') '2.80'
Set-TextValue $ws.Range('E42') '  -9.79%  '

# Row 43
Set-TextValue $ws.Range('D43') '0.0414'
Set-TextValue $ws.Range('E43') '  -7.62%  '

# Row 44
Set-TextValue $ws.Range('E44') '  -13.32%  '

# Row 45
Set-TextValue $ws.Range('D45') '2.705.75'
Set-TextValue $ws.Range('E45') '  -3.50%  '

# Row 46
Set-TextValue $ws.Range('E46') '  -4.56%  '

# Row 47
Set-TextValue $ws.Range('D47') '3.03'
Set-TextValue $ws.Range('E47') '  -6.89%  '

# Row 48
Set-TextValue $ws.Range('B48') 'WEMIXToken'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D48') '2.55'
Set-TextValue $ws.Range('E48') '  -7.39%  '

# Row 49
Set-TextValue $ws.Range('B49') 'Monero'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D49') '138.16'
Set-TextValue $ws.Range('E49') '  -3.00%  '

# Row 50
Set-TextValue $ws.Range('D50') '8.36'
Set-TextValue $ws.Range('E50') '  -11.36%  '

# Row 51
Set-TextValue $ws.Range('D51') '2.60'
Set-TextValue $ws.Range('E51') '  -7.47%  '
